$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-8 from serial 46066 (2026-02-13)
# to serial 46070 (2026-02-17), preserving existing date formatting/style.
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 46070
}
